$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A1:M1")
$r.Borders.LineStyle = 1
$r.Borders.Color = 16777215
